# Insert 4 new data rows (new rows 240-243) into the "Poroto granado" sheet.
# Everything that was previously at row 240 and below shifts down by 4 rows,
# which is handled automatically by Rows.Insert().

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at position 240 - existing row 240 (and everything
# below it) shifts down to make room, inheriting formatting from the
# row that used to be there (so column D keeps its date number format).
$ws.Rows("240:243").Insert()

# ---- New row 240 ----
$ws.Range("A240").Value = 6
$ws.Range("B240").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C240").Value = "Metropolitana"
$ws.Range("D240").Value = 44559
$ws.Range("E240").Value = 13
$ws.Range("F240").Value = 100112030
$ws.Range("G240").Value = "Poroto granado"
$ws.Range("H240").Value = "Sin especificar"
$ws.Range("I240").Value = "Primera"
$ws.Range("J240").Value = 1310
$ws.Range("K240").Value = 25000
$ws.Range("L240").Value = 27000
$ws.Range("M240").Value = 25855
$ws.Range("N240").Value = "`$/saco 25 kilos"
$ws.Range("O240").Value = "Región Metropolitana"
$ws.Range("P240").Value = 1034
$ws.Range("Q240").Value = 25
$ws.Range("R240").Value = "Hortaliza"

# ---- New row 241 ----
$ws.Range("A241").Value = 6
$ws.Range("B241").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C241").Value = "Metropolitana"
$ws.Range("D241").Value = 44559
$ws.Range("E241").Value = 13
$ws.Range("F241").Value = 100112030
$ws.Range("G241").Value = "Poroto granado"
$ws.Range("H241").Value = "Sin especificar"
$ws.Range("I241").Value = "Primera"
$ws.Range("J241").Value = 910
$ws.Range("K241").Value = 26000
$ws.Range("L241").Value = 28000
$ws.Range("M241").Value = 27231
$ws.Range("N241").Value = "`$/saco 25 kilos"
$ws.Range("O241").Value = "Región de O'Higgins"
$ws.Range("P241").Value = 1089
$ws.Range("Q241").Value = 25
$ws.Range("R241").Value = "Hortaliza"

# ---- New row 242 ----
$ws.Range("A242").Value = 6
$ws.Range("B242").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C242").Value = "Metropolitana"
$ws.Range("D242").Value = 44559
$ws.Range("E242").Value = 13
$ws.Range("F242").Value = 100112030
$ws.Range("G242").Value = "Poroto granado"
$ws.Range("H242").Value = "Sin especificar"
$ws.Range("I242").Value = "Segunda"
$ws.Range("J242").Value = 350
$ws.Range("K242").Value = 20000
$ws.Range("L242").Value = 20000
$ws.Range("M242").Value = 20000
$ws.Range("N242").Value = "`$/saco 25 kilos"
$ws.Range("O242").Value = "Región Metropolitana"
$ws.Range("P242").Value = 800
$ws.Range("Q242").Value = 25
$ws.Range("R242").Value = "Hortaliza"

# ---- New row 243 ----
$ws.Range("A243").Value = 6
$ws.Range("B243").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C243").Value = "Metropolitana"
$ws.Range("D243").Value = 44559
$ws.Range("E243").Value = 13
$ws.Range("F243").Value = 100112030
$ws.Range("G243").Value = "Poroto granado"
$ws.Range("H243").Value = "Sin especificar"
$ws.Range("I243").Value = "Segunda"
$ws.Range("J243").Value = 250
$ws.Range("K243").Value = 22000
$ws.Range("L243").Value = 22000
$ws.Range("M243").Value = 22000
$ws.Range("N243").Value = "`$/saco 25 kilos"
$ws.Range("O243").Value = "Región de O'Higgins"
$ws.Range("P243").Value = 880
$ws.Range("Q243").Value = 25
$ws.Range("R243").Value = "Hortaliza"
